$wb = $excel.ActiveWorkbook

# Rename the two worksheets (keeping their underlying sheetId / position).
$wb.Worksheets.Item(1).Name = "MEJORAR"
$wb.Worksheets.Item(2).Name = "PREMIUM"

# Make the second sheet ("PREMIUM") the active / selected tab, which moves
# tabSelected from sheet1 to sheet2 and sets workbookView activeTab="1".
$wb.Worksheets.Item(2).Activate()
